$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.160.29"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.763.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -3.04%  "

$ws.Range("E9").Value = "  -3.85%  "

$ws.Range("E10").Value = "  -7.34%  "

$ws.Range("E11").Value = "  +3.59%  "

$ws.Range("E12").Value = "  -2.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.252.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.808.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("E16").Value = "  -4.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.768.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.18%  "

$ws.Range("E19").Value = "  -3.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "361.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("E21").Value = "  -5.69%  "

$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.528"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.33%  "

$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0908"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.45%  "

$ws.Range("E30").Value = "  -4.25%  "

$ws.Range("E31").Value = "  +5.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.90%  "

$ws.Range("E33").Value = "  -3.55%  "

$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.08%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("E38").Value = "  -2.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "348.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "

$ws.Range("E41").Value = "  -2.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("E44").Value = "  -4.10%  "

$ws.Range("E45").Value = "  -3.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "137.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.41%  "

$ws.Range("E47").Value = "  -3.19%  "

$ws.Range("E48").Value = "  -3.35%  "

$ws.Range("E49").Value = "  -2.34%  "

$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "
